$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new "ARREARS" column label in I5 ---
$ws.Range("I5").Value = "ARREARS"

# --- Row 6 (MARTHA WAIRIMU / OK-1): add a 2000 B/F and a 2000 payment ---
$ws.Range("D6").Value = 2000
$ws.Range("G6").Value = 2000

# --- Row 7 (OMAE OMWENGA / OK-2): rent removed this period ---
$ws.Range("E7").ClearContents()

# --- Row 8 (RACHAEL GATHIRI / OK-3): add a 5000 B/F and a 5000 payment ---
$ws.Range("D8").Value = 5000
$ws.Range("G8").Value = 5000

# --- Row 9 (FAITH WANJIRU / OK-4): add a 1000 B/F, 6000 payment, 1000 arrears ---
$ws.Range("D9").Value = 1000
$ws.Range("G9").Value = 6000
$ws.Range("I9").Value = 1000

# --- Row 10 (OMAE OMWENGA / OK-5): full payment of 5000 ---
$ws.Range("G10").Value = 5000

# --- Row 12 (GEOFREY OMITI / OK-7): rent 6000, fully paid ---
$ws.Range("E12").Value = 6000
$ws.Range("G12").Value = 6000

# --- Row 13 (ANN MULI / OK-8): add 8000 B/F, 2500 arrears ---
$ws.Range("D13").Value = 8000
$ws.Range("I13").Value = 2500

# --- Row 20 totals: new ARREARS column total ---
$ws.Range("I20").Formula = "=SUM(I6:I19)"

# --- Summary section ---
$ws.Range("F26").Formula = "=H20"

$ws.Range("A28").Value = "ARREARS"

$ws.Range("A34").Value = "PAIN ON 14/12"
$ws.Range("C34").Value = 33400
$ws.Range("E34").Value = "PAIN ON 14/12"
$ws.Range("G34").Value = 33400

# --- View state: scroll + selection as left by the editor ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G7").Select()
